$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($addr, $text, $refAddr) {
    $ws.Range($addr).Value = $text
    $ws.Range($refAddr).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# --- Step 1: drop the two trailing rows so the sheet shrinks from 25 to 23 rows ---
$ws.Rows.Item(25).Delete()
$ws.Rows.Item(24).Delete()

# --- Step 2: row 10 ("Objetivos:") body text is replaced by a professor id/name ---
Set-CellText "B10" "4780627 - Ana Lucia Gabas Ferreira" "B3"
Set-CellText "C10" "4780627 - Ana Lucia Gabas Ferreira" "C3"

# --- Step 3: rebuild rows 13-23 content/structure ---
# Row 13
Set-CellText "A13" "Programa resumido:" "A3"
Set-CellText "B13" "'01/01/2022" "B3"
Set-CellText "C13" "'01/01/2022" "C3"
$ws.Rows.Item(13).RowHeight = 60

# Row 14
Set-CellText "A14" "Short syllabus:" "A3"
Set-CellText "B14" "The student must develop the Graduate Work 2 following the model chosen and already used in the Graduate Work 1. When choosing the model, the student must develop scientific content or applicable product, both related to Environmental Engineering, and at the end, the Graduate Work must be evaluated by an evaluating panel." "B3"
Set-CellText "C14" "The student must develop the Graduate Work 2 following the model chosen and already used in the Graduate Work 1. When choosing the model, the student must develop scientific content or applicable product, both related to Environmental Engineering, and at the end, the Graduate Work must be evaluated by an evaluating panel." "C3"
$ws.Rows.Item(14).RowHeight = 60

# Row 15
Set-CellText "A15" "Programa:" "A3"
Set-CellText "B15" "4780627 - Ana Lucia Gabas Ferreira" "B3"
Set-CellText "C15" "4780627 - Ana Lucia Gabas Ferreira" "C3"
$ws.Rows.Item(15).RowHeight = 120

# Row 16
Set-CellText "A16" "Syllabus:" "A3"
Set-CellText "B16" "For both models (article or product): The student must continue with the development of the Work initiated in the Graduate Work 1. The course program consists of the following steps: 1) Development of the theme based on the schedule activities approved in the Graduate Work 1. 2) Development of the final text, according to the model provided by those responsible for the discipline. 3) Delivery of the final version of the text, with the approval of the advisor and with the indication of the evaluation board. 4) Evaluation and grading by the examining board, which may or may not be the same previously composed in Graduate Work 1, at the discretion of the advisor." "B3"
Set-CellText "C16" "For both models (article or product): The student must continue with the development of the Work initiated in the Graduate Work 1. The course program consists of the following steps: 1) Development of the theme based on the schedule activities approved in the Graduate Work 1. 2) Development of the final text, according to the model provided by those responsible for the discipline. 3) Delivery of the final version of the text, with the approval of the advisor and with the indication of the evaluation board. 4) Evaluation and grading by the examining board, which may or may not be the same previously composed in Graduate Work 1, at the discretion of the advisor." "C3"
$ws.Rows.Item(16).RowHeight = 120

# Row 17
Set-CellText "A17" "Avaliação:" "A3"
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()
$ws.Rows.Item(17).EntireRow.AutoFit()

# Row 18
Set-CellText "A18" "Método:" "A3"
Set-CellText "B18" "7455355 - Robson da Silva Rocha" "B3"
Set-CellText "C18" "7455355 - Robson da Silva Rocha" "C3"
$ws.Rows.Item(18).RowHeight = 60

# Row 19
Set-CellText "A19" "Critério:" "A3"
Set-CellText "B19" "Para ambos os modelos (artigo ou produto): O aluno deverá apresentar um artigo científico seguindo modelo fornecido pelos responsáveis da disciplina. A versão final do texto deverá ser aprovada pelo orientador no ato da submissão e deverá ser apresentado para banca de avaliação composta por dois doutores ou especialistas na área do projeto. Alternativamente, pode-se apresentar o texto com formatação de outra revista desde que seja anexado o comprovante de submissão do artigo, nesse caso, o texto apresentado pode seguir as regras de formatação da revista escolhida. A revista escolhida deve ser indexada por algum sistema de base de dados (Web of Science, Scopus, Scielo, Cinahl, Medline, etc). Alternativamente, artigos aprovados em revistas da área de estudo, até a semana anterior a apresentação, desobrigam o aluno a apresentar o trabalho para a banca e nesse caso, a entrega do artigo deve ser acompanhada pela comprovação do aceite do trabalho. A revista escolhida deve ser indexada por algum sistema de base de dados (Web of Science, Scopus, Scielo, Cinahl, Medline, etc)" "B3"
Set-CellText "C19" "Para ambos os modelos (artigo ou produto): O aluno deverá apresentar um artigo científico seguindo modelo fornecido pelos responsáveis da disciplina. A versão final do texto deverá ser aprovada pelo orientador no ato da submissão e deverá ser apresentado para banca de avaliação composta por dois doutores ou especialistas na área do projeto. Alternativamente, pode-se apresentar o texto com formatação de outra revista desde que seja anexado o comprovante de submissão do artigo, nesse caso, o texto apresentado pode seguir as regras de formatação da revista escolhida. A revista escolhida deve ser indexada por algum sistema de base de dados (Web of Science, Scopus, Scielo, Cinahl, Medline, etc). Alternativamente, artigos aprovados em revistas da área de estudo, até a semana anterior a apresentação, desobrigam o aluno a apresentar o trabalho para a banca e nesse caso, a entrega do artigo deve ser acompanhada pela comprovação do aceite do trabalho. A revista escolhida deve ser indexada por algum sistema de base de dados (Web of Science, Scopus, Scielo, Cinahl, Medline, etc)" "C3"
$ws.Rows.Item(19).RowHeight = 60

# Row 20
Set-CellText "A20" "Norma de recuperação:" "A3"
Set-CellText "B20" "Avaliação e emissão de parecer pela banca avaliadora e pelo orientador, com atribuição de nota única final.Fica sob responsabilidade do orientador a verificação de ocorrência de plágio utilizando software apropriado e avaliação em Comitê de Ética, quando exigido, via cadastro na Plataforma Brasil." "B3"
Set-CellText "C20" "Avaliação e emissão de parecer pela banca avaliadora e pelo orientador, com atribuição de nota única final.Fica sob responsabilidade do orientador a verificação de ocorrência de plágio utilizando software apropriado e avaliação em Comitê de Ética, quando exigido, via cadastro na Plataforma Brasil." "C3"
$ws.Rows.Item(20).RowHeight = 60

# Row 21
Set-CellText "A21" "Bibliografia:" "A3"
Set-CellText "B21" "Não há." "B3"
Set-CellText "C21" "Não há." "C3"
$ws.Rows.Item(21).RowHeight = 120

# Row 22
Set-CellText "A22" "Requisitos:" "A3"
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Rows.Item(22).EntireRow.AutoFit()

# Row 23
Set-CellText "B23" "LOB1223 -  Trabalho de Graduação I  (Requisito fraco)
" "B3"
Set-CellText "C23" "LOB1223 -  Trabalho de Graduação I  (Requisito fraco)
" "C3"
$ws.Range("A23").Clear()
$ws.Rows.Item(23).RowHeight = 30
